$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Columns("A").Insert()

# Populate the new column A with Runmode header + Yes/Yes/No/No values
$ws.Range("A1").Value = "Runmode"
$ws.Range("A2").Value = "Yes"
$ws.Range("A3").Value = "Yes"
$ws.Range("A4").Value = "No"
$ws.Range("A5").Value = "No"

# Adjust column widths for the new columns A and B
$ws.Range("A1").ColumnWidth = 17.75
$ws.Range("B1").ColumnWidth = 22.75

# Reset selection to the default top-left cell
$ws.Range("A1").Select() | Out-Null

